# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that a leading "System" token is moved to the end of the comma-separated
# list instead of being listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val

    if ($text.StartsWith("System, ")) {
        $parts = $text -split ", "
        $first = $parts[0]
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + @($first)
        $newText = $newParts -join ", "
        $cell.Value = $newText
    }
}
